# Update column G ("K") values per regenerated save_data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 0
    6 = 1
    7 = 2
    8 = 2
    9 = 2
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 2
    19 = 0
    20 = 0
    21 = 2
    22 = 2
    23 = 2
    24 = 1
    25 = 2
    26 = 0
    27 = 0
    28 = 2
    30 = 1
    31 = 0
    32 = 1
    33 = 2
    34 = 2
    35 = 0
    36 = 0
    37 = 1
    38 = 1
    39 = 1
    40 = 2
    41 = 1
    42 = 2
    43 = 0
    44 = 1
    45 = 0
    46 = 0
    47 = 1
    48 = 2
    49 = 0
    50 = 1
    51 = 1
    52 = 1
    53 = 2
    54 = 1
    55 = 3
    56 = 0
    57 = 1
    58 = 0
    59 = 1
    60 = 1
    61 = 0
    62 = 2
    63 = 0
    64 = 1
    65 = 1
    66 = 1
    67 = 0
    68 = 0
    69 = 1
    70 = 1
    71 = 1
    72 = 1
    76 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
